$wb = $excel.ActiveWorkbook

# The two sheets that receive the new "UC" benefit columns/rows
$sheetNames = @("UK_DHE_MCS2_Males", "UK_DHE_MCS2_Females")

# New regressor names -> go in column A of new rows 10-14, and as new
# column headers K1:O1 (shared strings 42-46)
$newNames = @(
    "D_Econ_benefits_UC_Lhw_ZERO",
    "D_Econ_benefits_UC_Lhw_TEN",
    "D_Econ_benefits_UC_Lhw_TWENTY",
    "D_Econ_benefits_UC_Lhw_THIRTY",
    "D_Econ_benefits_UC_Lhw_FORTY"
)

# Column letters for the 5 new columns (K..O)
$newCols = @("K", "L", "M", "N", "O")

# Constant (column B) for each new row 10-14
$constants = @(-2.2400000000000002, -1.27, -1.0740000000000001, -0.55900000000000005, -0.76900000000000002)

# Diagonal coefficient for each new row/col pair (row 10 <-> K, row 11 <-> L, ...)
# (written in plain decimal since the PowerShell parser here doesn't accept
# scientific-notation numeric literals)
$diag = @(0.090920710120782997, 0.40933881715951598, 0.23123340535193601, 0.22708142440649701, 0.140625)

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # --- Header row: K1:O1 ---
    for ($i = 0; $i -lt 5; $i++) {
        $ws.Range($newCols[$i] + "1").Value = $newNames[$i]
    }

    # --- Existing rows 2-9: fill new K:O columns with 0 ---
    for ($r = 2; $r -le 9; $r++) {
        for ($i = 0; $i -lt 5; $i++) {
            $ws.Range($newCols[$i] + $r).Value = 0
        }
    }

    # --- New rows 10-14 ---
    for ($i = 0; $i -lt 5; $i++) {
        $r = 10 + $i
        $ws.Range("A" + $r).Value = $newNames[$i]
        $ws.Range("B" + $r).Value = $constants[$i]
        for ($c = 3; $c -le 10; $c++) {
            $ws.Cells.Item($r, $c).Value = 0
        }
        for ($j = 0; $j -lt 5; $j++) {
            if ($j -eq $i) {
                $ws.Range($newCols[$j] + $r).Value = $diag[$i]
            } else {
                $ws.Range($newCols[$j] + $r).Value = 0
            }
        }
    }
}

# --- View / selection state ---
$wsMales = $wb.Worksheets.Item("UK_DHE_MCS2_Males")
$wsMales.Range("K1:O9").Select()

$wsFemales = $wb.Worksheets.Item("UK_DHE_MCS2_Females")
$wsFemales.Activate()
$wsFemales.Range("O4").Select()
